# Apply cell updates per the crypto price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    # Force the cell to stay a text string even when the text looks numeric
    # (e.g. "582.56"), matching the inlineStr cells produced by the source data,
    # then restore the default (unstyled) cell formatting so no style is changed.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "66.209.40"
$ws.Range("E2").Value = "  +6.57%  "
$ws.Range("D3").Value = "3.008.72"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "582.56"
$ws.Range("E5").Value = "  +2.82%  "
Set-TextValue $ws.Range("D6") "163.07"
$ws.Range("E6").Value = "  +13.32%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.003.43"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D9") "0.517"
$ws.Range("E9").Value = "  +3.48%  "
Set-TextValue $ws.Range("D10") "6.69"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  +6.55%  "
Set-TextValue $ws.Range("D14") "34.59"
$ws.Range("E14").Value = "  +6.41%  "
$ws.Range("D16").Value = "66.215.76"
$ws.Range("E16").Value = "  +6.69%  "
$ws.Range("D17").Value = "3.509.18"
$ws.Range("E17").Value = "  +3.63%  "
Set-TextValue $ws.Range("D18") "6.93"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("D19").Value = "3.013.25"
$ws.Range("E19").Value = "  +4.00%  "
Set-TextValue $ws.Range("D20") "454.46"
$ws.Range("E20").Value = "  +5.80%  "
Set-TextValue $ws.Range("D21") "13.89"
$ws.Range("E21").Value = "  +6.70%  "
Set-TextValue $ws.Range("D22") "0.687"
$ws.Range("E22").Value = "  +4.67%  "
$ws.Range("E23").Value = "  +7.28%  "
Set-TextValue $ws.Range("D24") "82.28"
$ws.Range("E24").Value = "  +4.72%  "
$ws.Range("E25").Value = "  +14.61%  "
$ws.Range("E26").Value = "  +2.91%  "
$ws.Range("E27").Value = "  +4.91%  "
$ws.Range("E28").Value = "  +0.02%  "
Set-TextValue $ws.Range("D29") "8.15"
$ws.Range("E29").Value = "  +17.68%  "
$ws.Range("E30").Value = "  +19.40%  "
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("E32").Value = "  +4.07%  "
Set-TextValue $ws.Range("D33") "27.20"
$ws.Range("E33").Value = "  +6.04%  "
$ws.Range("E34").Value = "  +5.03%  "
Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.08%  "
Set-TextValue $ws.Range("D36") "0.992"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("E37").Value = "  +7.82%  "
$ws.Range("E38").Value = "  +15.84%  "
Set-TextValue $ws.Range("D39") "3.01"
$ws.Range("E39").Value = "  +2.38%  "
Set-TextValue $ws.Range("D40") "49.92"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("E41").Value = "  +16.26%  "
$ws.Range("E42").Value = "  +8.18%  "
Set-TextValue $ws.Range("D43") "43.99"
$ws.Range("E43").Value = "  +8.11%  "
Set-TextValue $ws.Range("D44") "8.44"
$ws.Range("E44").Value = "  +4.01%  "
Set-TextValue $ws.Range("D45") "396.51"
$ws.Range("E45").Value = "  +14.72%  "
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("D47").Value = "2.795.45"
$ws.Range("E47").Value = "  +2.98%  "
Set-TextValue $ws.Range("D48") "133.97"
$ws.Range("E48").Value = "  +0.47%  "
Set-TextValue $ws.Range("D50") "23.99"
$ws.Range("E50").Value = "  +12.12%  "
$ws.Range("E51").Value = "  +4.37%  "
